$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132; all existing rows 132..277 shift down to 133..278
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new weekly data point
$ws.Range("A132").Value = 5
$ws.Range("B132").Value = "Macroferia Regional de Talca"
$ws.Range("C132").Value = "Maule"
$ws.Range("D132").Value = 44629
$ws.Range("E132").Value = 7
$ws.Range("F132").Value = 100114014
$ws.Range("G132").Value = "Betarraga"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 800
$ws.Range("L132").Value = 800
$ws.Range("M132").Value = 800
$ws.Range("N132").Value = "$/paquete 5 unidades"
$ws.Range("O132").Value = "Región del Maule"
$ws.Range("P132").Value = 160
$ws.Range("Q132").Value = 5
$ws.Range("R132").Value = "Hortaliza"
